$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header separator row ("//") into the two new columns (F1:G1)
$ws.Range("F1").Value = "//"
$ws.Range("G1").Value = "//"

# Replace the visit record in row 2 with the new data
$ws.Range("D2").Value = "19.09.2022"
$ws.Range("D2").ClearFormats()

$ws.Range("E2:G2").NumberFormat = "@"
$ws.Range("E2").Value = "19.22"
$ws.Range("F2").Value = "wiktor.k.2002@icloud.com"
$ws.Range("G2").Value = "01"
$ws.Range("E2:G2").ClearFormats()

[void]$ws.Range("G1").Select()
